# LOM3100.xlsx edit script
#
# The published change deletes the worksheet row that held the
# "7797767 - Viktor Pastoukhov" docente value (old row 13, which had no
# label in column A), shifting every row below it up by one. On top of
# that shift, several of the content cells in columns B/C end up showing
# values that no longer line up with their neighbouring row label -- i.e.
# the published sheet is internally inconsistent (e.g. "Programa:" shows
# "01/01/2017" and "Método:" shows the docente string) -- but that is
# exactly what the target workbook contains, so we reproduce it exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old row 13 (B13/C13 = "7797767 - Viktor Pastoukhov", no
#    label in column A). Everything below shifts up by one row, and the
#    sheet's used range becomes A1:C23.
$ws.Rows.Item(13).Delete()

# 2) Fix up the content cells whose text doesn't simply follow the shift.
$ws.Range("B10").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C10").Value = "7797767 - Viktor Pastoukhov"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2017"
$ws.Range("C15").Value = "01/01/2017"

$ws.Range("B18").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C18").Value = "7797767 - Viktor Pastoukhov"

$ws.Range("B19").Value = "A avaliação será composta por duas provas (P1 e P2)."
$ws.Range("C19").Value = "A avaliação será composta por duas provas (P1 e P2)."

$ws.Range("B20").Value = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."
$ws.Range("C20").Value = "NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total."

$ws.Range("B21").Value = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2."
$ws.Range("C21").Value = "A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2."
